$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert two new columns before column I (9), shifting old I..L -> K..N
$ws.Columns("I:J").Insert()

# 2. Insert a new row before row 16 (the numbered sub-header row),
#    shifting old row16(data) -> 17, row19 -> 20, row24 -> 25
$ws.Rows("16:16").Insert()

# 3. New header labels for the two inserted columns (row 15 already
#    picked up the correct bordered/centered style from its left neighbour)
$ws.Range("I15").Value = "Số thửa"
$ws.Range("J15").Value = "Tờ bản đồ"

# 4. Build the new numbered sub-header row 16 (A16:N16 = 1..14),
#    matching the bordered/centered style used by the row-15 headers.
$ws.Range("A15").Copy()
$ws.Range("A16:N16").PasteSpecial(-4122)
for ($i = 1; $i -le 14; $i++) {
    $ws.Cells.Item(16, $i).Value = $i
}
$ws.Rows("16:16").RowHeight = 16.2

# 5. Turn the new row into the AutoFilter header row.
$ws.Range("A16:N16").AutoFilter(1)

# 6. Register the (hidden) _FilterDatabase defined name that Excel
#    creates for the AutoFilter range, scoped to this sheet.
$n = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$16:`$N`$16")
$n.Visible = $false

# 7. Restore/update the active selection shown in the sheet view.
$ws.Range("A16:N16").Select()
